$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.801.41"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "'2.318.51"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'97.20"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").Value = "'271.85"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.625"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'45.36"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "'2.658.10"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "'15.52"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "'0.882"
$ws.Range("E16").Value = "  +10.00%  "
$ws.Range("D17").Value = "'2.325.14"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "'43.765.28"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = "  +3.76%  "
$ws.Range("D20").Value = "'6.39"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("D21").Value = "'73.30"
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("D22").Value = "'240.25"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").Value = "'9.41"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'3.50"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'38.09"
$ws.Range("E30").Value = "  -7.99%  "
$ws.Range("E31").Value = "  +6.85%  "
$ws.Range("D32").Value = "'174.72"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").Value = "'0.0910"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "'5.47"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  -5.50%  "
$ws.Range("E40").Value = "  +8.94%  "
$ws.Range("D41").Value = "'2.37"
$ws.Range("E41").Value = "  +8.59%  "
$ws.Range("E42").Value = "  +19.34%  "
$ws.Range("E43").Value = "  -5.51%  "
$ws.Range("E44").Value = "  +10.01%  "
$ws.Range("D45").Value = "'62.43"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'100.33"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("E50").Value = "  +14.94%  "
$ws.Range("D51").Value = "'2.545.47"
$ws.Range("E51").Value = "  +3.91%  "
